$d = $word.ActiveDocument

$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function New-ParaXml([string]$body) {
    return $xmlHeader + $body + $xmlFooter
}

# ---------------------------------------------------------------------------
# 1) "DetectedActivitiesIntentService.java:" paragraph -- merge the two bold
#    runs ("DetectedActivitiesIntentService" + ".java") into a single run.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "DetectedActivitiesIntentService.java:*") {
        $target = $d.Paragraphs.Item($i)
        break
    }
}
$body = '<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>DetectedActivitiesIntentService.java</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t xml:space="preserve">Intent service to handle incoming intents generated from the Google Activity Recognition API. </w:t></w:r></w:p>'
$null = $target.Range.InsertXML((New-ParaXml $body))

# ---------------------------------------------------------------------------
# 2) "GeofenceErrorMessages.java:" paragraph -- drop the empty-paragraph
#    <w:pPr><w:rPr><w:b/></w:rPr></w:pPr> mark formatting and append the new
#    descriptive sentence (with a spell-checked "Geofence" run).
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "GeofenceErrorMessages.java:*") {
        $target = $d.Paragraphs.Item($i)
        break
    }
}
$body = '<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>GeofenceErrorMessages.java</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> Helper class for debugging errors associated with </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Geofence</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> creation.</w:t></w:r></w:p>'
$null = $target.Range.InsertXML((New-ParaXml $body))

# ---------------------------------------------------------------------------
# 3) "GeofenceTransitionsIntentService.java:" paragraph -- drop the mark
#    formatting and append the descriptive sentence (two spell-checked runs:
#    "geofence" and "TrialActivity").
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "GeofenceTransitionsIntentService.java:*") {
        $target = $d.Paragraphs.Item($i)
        break
    }
}
$body = '<w:p><w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>GeofenceTransitionsIntentService.java:</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Helper class to act as a listener for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>geofence</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> transition changes. Works with the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TrialActivity</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to determine when to stop recording. </w:t></w:r></w:p>'
$null = $target.Range.InsertXML((New-ParaXml $body))

# ---------------------------------------------------------------------------
# 4) "MapsActivity.java:" paragraph -- merge the bold runs into one, drop the
#    mark formatting, and append the descriptive sentence.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "MapsActivity*java:*") {
        $target = $d.Paragraphs.Item($i)
        break
    }
}
$body = '<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>MapsActivity.java:</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Main activity that is shown to the user. Shows parking lots as markers on a Google Map. The color of the marker determines the availability. Clicking on a parking lot provides additional information relevant to the parking status. </w:t></w:r></w:p>'
$null = $target.Range.InsertXML((New-ParaXml $body))

# ---------------------------------------------------------------------------
# 5) "ParkingAnalyzer.java:" paragraph -- merge the bold runs into one, drop
#    the mark formatting, and append the descriptive sentence.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "ParkingAnalyzer*java:*") {
        $target = $d.Paragraphs.Item($i)
        break
    }
}
$body = '<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>ParkingAnalyzer.java:</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Primary class that handles all the estimation of parking availability. Contains all the feature identification and classification. </w:t></w:r></w:p>'
$null = $target.Range.InsertXML((New-ParaXml $body))

# ---------------------------------------------------------------------------
# 6) "ParkingParser.java:" and "TrialActivity.java:" paragraphs -- these are
#    the last two paragraphs in the body, so they are replaced together in a
#    single InsertXML call (spanning both) to avoid leaving a stray empty
#    paragraph behind (InsertXML keeps the body's final paragraph mark when
#    it is included inside the replaced range). The _GoBack bookmark, which
#    used to sit on the "GeofenceErrorMessages.java:" paragraph, now moves to
#    the end of the final ("TrialActivity.java:") paragraph.
# ---------------------------------------------------------------------------
$parserPara = $null
$trialPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "ParkingParser*java:*") {
        $parserPara = $d.Paragraphs.Item($i)
    }
    if ($d.Paragraphs.Item($i).Range.Text -like "TrialActivity*java:*") {
        $trialPara = $d.Paragraphs.Item($i)
    }
}
$combinedRange = $d.Range($parserPara.Range.Start, $trialPara.Range.End)
$body = '<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>ParkingParser.java:</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Helper class used in ParkingAnalyzer.java. Used to take read csv files into usable arrays that can be interpreted in ParkingAnalyzer.java. </w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">TrialActivity.java: </w:t></w:r><w:r><w:t xml:space="preserve">Class that was used in testing and collection of data. Saves sensor data into csv files for later interpretation and analysis. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$null = $combinedRange.InsertXML((New-ParaXml $body))
# NB: the old _GoBack bookmark that sat on the "GeofenceErrorMessages.java:"
# paragraph was already removed above when that paragraph's range was
# replaced wholesale (step 2), so there is nothing left to clean up here.

Write-Output "Done"
